{"js": "// Replace each two-digit multiplication \"expression=result\" cell text in the\n// table with its updated value. The mapping below is old text -> new text,\n// taken in document order (one per table cell, excluding the date heading).\nconst pairs = [\n  [\"52\u00d712=624\", \"61\u00d727=1647\"],\n  [\"78\u00d723=1794\", \"99\u00d718=1782\"],\n  [\"84\u00d766=5544\", \"68\u00d793=6324\"],\n  [\"31\u00d780=2480\", \"91\u00d729=2639\"],\n  [\"47\u00d726=1222\", \"90\u00d782=7380\"],\n  [\"21\u00d766=1386\", \"100\u00d789=8900\"],\n  [\"11\u00d741=451\", \"31\u00d724=744\"],\n  [\"24\u00d735=840\", \"40\u00d736=1440\"],\n  [\"80\u00d789=7120\", \"88\u00d726=2288\"],\n  [\"80\u00d747=3760\", \"13\u00d732=416\"],\n  [\"36\u00d761=2196\", \"33\u00d718=594\"],\n  [\"63\u00d748=3024\", \"20\u00d723=460\"],\n  [\"43\u00d782=3526\", \"12\u00d799=1188\"],\n  [\"44\u00d740=1760\", \"46\u00d769=3174\"],\n  [\"56\u00d755=3080\", \"13\u00d736=468\"],\n  [\"17\u00d781=1377\", \"51\u00d766=3366\"],\n  [\"98\u00d794=9212\", \"87\u00d792=8004\"],\n  [\"22\u00d729=638\", \"27\u00d766=1782\"],\n  [\"19\u00d748=912\", \"45\u00d722=990\"],\n  [\"22\u00d744=968\", \"22\u00d778=1716\"],\n  [\"18\u00d724=432\", \"84\u00d795=7980\"],\n  [\"91\u00d774=6734\", \"22\u00d791=2002\"],\n  [\"57\u00d755=3135\", \"58\u00d727=1566\"],\n  [\"67\u00d711=737\", \"82\u00d771=5822\"],\n  [\"25\u00d725=625\", \"38\u00d787=3306\"],\n  [\"41\u00d763=2583\", \"34\u00d780=2720\"],\n  [\"73\u00d738=2774\", \"94\u00d719=1786\"],\n  [\"50\u00d777=3850\", \"96\u00d737=3552\"],\n  [\"35\u00d720=700\", \"35\u00d747=1645\"],\n  [\"97\u00d766=6402\", \"67\u00d717=1139\"],\n  [\"33\u00d799=3267\", \"74\u00d768=5032\"],\n  [\"27\u00d745=1215\", \"20\u00d766=1320\"],\n  [\"24\u00d712=288\", \"38\u00d749=1862\"],\n  [\"95\u00d762=5890\", \"82\u00d751=4182\"],\n  [\"43\u00d765=2795\", \"100\u00d742=4200\"],\n  [\"27\u00d741=1107\", \"49\u00d786=4214\"],\n  [\"54\u00d737=1998\", \"85\u00d725=2125\"],\n  [\"13\u00d761=793\", \"10\u00d710=100\"],\n  [\"55\u00d746=2530\", \"65\u00d749=3185\"],\n  [\"67\u00d785=5695\", \"95\u00d781=7695\"],\n  [\"41\u00d789=3649\", \"28\u00d737=1036\"],\n  [\"78\u00d779=6162\", \"43\u00d760=2580\"],\n  [\"91\u00d757=5187\", \"22\u00d731=682\"],\n  [\"25\u00d799=2475\", \"33\u00d797=3201\"],\n  [\"47\u00d769=3243\", \"97\u00d799=9603\"],\n  [\"97\u00d725=2425\", \"31\u00d768=2108\"],\n  [\"78\u00d793=7254\", \"51\u00d799=5049\"],\n  [\"71\u00d789=6319\", \"100\u00d757=5700\"],\n  [\"95\u00d744=4180\", \"20\u00d756=1120\"],\n  [\"10\u00d738=380\", \"38\u00d731=1178\"],\n  [\"84\u00d745=3780\", \"74\u00d711=814\"],\n  [\"85\u00d747=3995\", \"82\u00d733=2706\"],\n  [\"78\u00d736=2808\", \"19\u00d717=323\"],\n  [\"93\u00d789=8277\", \"29\u00d782=2378\"],\n  [\"43\u00d757=2451\", \"26\u00d723=598\"],\n  [\"38\u00d748=1824\", \"42\u00d795=3990\"],\n  [\"71\u00d726=1846\", \"81\u00d772=5832\"],\n  [\"53\u00d757=3021\", \"95\u00d753=5035\"],\n  [\"32\u00d744=1408\", \"96\u00d710=960\"],\n  [\"51\u00d769=3519\", \"15\u00d770=1050\"],\n  [\"57\u00d747=2679\", \"80\u00d728=2240\"],\n  [\"56\u00d765=3640\", \"54\u00d756=3024\"],\n  [\"93\u00d767=6231\", \"82\u00d771=5822\"],\n  [\"72\u00d748=3456\", \"54\u00d719=1026\"],\n  [\"68\u00d758=3944\", \"55\u00d715=825\"],\n  [\"14\u00d797=1358\", \"55\u00d760=3300\"],\n  [\"17\u00d720=340\", \"50\u00d755=2750\"],\n  [\"35\u00d737=1295\", \"68\u00d775=5100\"],\n  [\"69\u00d733=2277\", \"57\u00d781=4617\"],\n  [\"33\u00d777=2541\", \"66\u00d767=4422\"],\n  [\"90\u00d799=8910\", \"16\u00d794=1504\"],\n  [\"65\u00d768=4420\", \"83\u00d785=7055\"],\n  [\"78\u00d762=4836\", \"51\u00d742=2142\"],\n  [\"44\u00d786=3784\", \"40\u00d789=3560\"],\n  [\"11\u00d724=264\", \"25\u00d721=525\"],\n  [\"50\u00d730=1500\", \"33\u00d714=462\"],\n  [\"60\u00d738=2280\", \"30\u00d796=2880\"],\n  [\"68\u00d722=1496\", \"90\u00d7100=9000\"],\n  [\"96\u00d733=3168\", \"71\u00d767=4757\"],\n  [\"13\u00d729=377\", \"79\u00d745=3555\"],\n  [\"22\u00d724=528\", \"28\u00d777=2156\"],\n  [\"53\u00d737=1961\", \"19\u00d741=779\"],\n  [\"38\u00d737=1406\", \"38\u00d755=2090\"],\n  [\"95\u00d782=7790\", \"30\u00d763=1890\"],\n  [\"76\u00d762=4712\", \"41\u00d7100=4100\"],\n  [\"78\u00d713=1014\", \"59\u00d780=4720\"],\n  [\"73\u00d789=6497\", \"71\u00d767=4757\"],\n  [\"60\u00d793=5580\", \"30\u00d761=1830\"],\n  [\"12\u00d777=924\", \"54\u00d793=5022\"],\n  [\"37\u00d741=1517\", \"36\u00d789=3204\"],\n  [\"88\u00d785=7480\", \"61\u00d799=6039\"],\n  [\"28\u00d794=2632\", \"98\u00d713=1274\"],\n  [\"64\u00d797=6208\", \"37\u00d738=1406\"],\n  [\"72\u00d719=1368\", \"38\u00d780=3040\"],\n  [\"41\u00d743=1763\", \"86\u00d720=1720\"],\n  [\"50\u00d786=4300\", \"47\u00d795=4465\"],\n  [\"50\u00d793=4650\", \"67\u00d790=6030\"],\n  [\"19\u00d710=190\", \"82\u00d757=4674\"],\n  [\"84\u00d791=7644\", \"36\u00d715=540\"],\n  [\"89\u00d797=8633\", \"41\u00d777=3157\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace each two-digit multiplication \"expression=result\" cell text in the\n# table with its updated value. The mapping below is old text -> new text,\n# taken in document order (one per table cell, excluding the date heading).\n$pairs = @(\n    @(\"52\u00d712=624\", \"61\u00d727=1647\"),\n    @(\"78\u00d723=1794\", \"99\u00d718=1782\"),\n    @(\"84\u00d766=5544\", \"68\u00d793=6324\"),\n    @(\"31\u00d780=2480\", \"91\u00d729=2639\"),\n    @(\"47\u00d726=1222\", \"90\u00d782=7380\"),\n    @(\"21\u00d766=1386\", \"100\u00d789=8900\"),\n    @(\"11\u00d741=451\", \"31\u00d724=744\"),\n    @(\"24\u00d735=840\", \"40\u00d736=1440\"),\n    @(\"80\u00d789=7120\", \"88\u00d726=2288\"),\n    @(\"80\u00d747=3760\", \"13\u00d732=416\"),\n    @(\"36\u00d761=2196\", \"33\u00d718=594\"),\n    @(\"63\u00d748=3024\", \"20\u00d723=460\"),\n    @(\"43\u00d782=3526\", \"12\u00d799=1188\"),\n    @(\"44\u00d740=1760\", \"46\u00d769=3174\"),\n    @(\"56\u00d755=3080\", \"13\u00d736=468\"),\n    @(\"17\u00d781=1377\", \"51\u00d766=3366\"),\n    @(\"98\u00d794=9212\", \"87\u00d792=8004\"),\n    @(\"22\u00d729=638\", \"27\u00d766=1782\"),\n    @(\"19\u00d748=912\", \"45\u00d722=990\"),\n    @(\"22\u00d744=968\", \"22\u00d778=1716\"),\n    @(\"18\u00d724=432\", \"84\u00d795=7980\"),\n    @(\"91\u00d774=6734\", \"22\u00d791=2002\"),\n    @(\"57\u00d755=3135\", \"58\u00d727=1566\"),\n    @(\"67\u00d711=737\", \"82\u00d771=5822\"),\n    @(\"25\u00d725=625\", \"38\u00d787=3306\"),\n    @(\"41\u00d763=2583\", \"34\u00d780=2720\"),\n    @(\"73\u00d738=2774\", \"94\u00d719=1786\"),\n    @(\"50\u00d777=3850\", \"96\u00d737=3552\"),\n    @(\"35\u00d720=700\", \"35\u00d747=1645\"),\n    @(\"97\u00d766=6402\", \"67\u00d717=1139\"),\n    @(\"33\u00d799=3267\", \"74\u00d768=5032\"),\n    @(\"27\u00d745=1215\", \"20\u00d766=1320\"),\n    @(\"24\u00d712=288\", \"38\u00d749=1862\"),\n    @(\"95\u00d762=5890\", \"82\u00d751=4182\"),\n    @(\"43\u00d765=2795\", \"100\u00d742=4200\"),\n    @(\"27\u00d741=1107\", \"49\u00d786=4214\"),\n    @(\"54\u00d737=1998\", \"85\u00d725=2125\"),\n    @(\"13\u00d761=793\", \"10\u00d710=100\"),\n    @(\"55\u00d746=2530\", \"65\u00d749=3185\"),\n    @(\"67\u00d785=5695\", \"95\u00d781=7695\"),\n    @(\"41\u00d789=3649\", \"28\u00d737=1036\"),\n    @(\"78\u00d779=6162\", \"43\u00d760=2580\"),\n    @(\"91\u00d757=5187\", \"22\u00d731=682\"),\n    @(\"25\u00d799=2475\", \"33\u00d797=3201\"),\n    @(\"47\u00d769=3243\", \"97\u00d799=9603\"),\n    @(\"97\u00d725=2425\", \"31\u00d768=2108\"),\n    @(\"78\u00d793=7254\", \"51\u00d799=5049\"),\n    @(\"71\u00d789=6319\", \"100\u00d757=5700\"),\n    @(\"95\u00d744=4180\", \"20\u00d756=1120\"),\n    @(\"10\u00d738=380\", \"38\u00d731=1178\"),\n    @(\"84\u00d745=3780\", \"74\u00d711=814\"),\n    @(\"85\u00d747=3995\", \"82\u00d733=2706\"),\n    @(\"78\u00d736=2808\", \"19\u00d717=323\"),\n    @(\"93\u00d789=8277\", \"29\u00d782=2378\"),\n    @(\"43\u00d757=2451\", \"26\u00d723=598\"),\n    @(\"38\u00d748=1824\", \"42\u00d795=3990\"),\n    @(\"71\u00d726=1846\", \"81\u00d772=5832\"),\n    @(\"53\u00d757=3021\", \"95\u00d753=5035\"),\n    @(\"32\u00d744=1408\", \"96\u00d710=960\"),\n    @(\"51\u00d769=3519\", \"15\u00d770=1050\"),\n    @(\"57\u00d747=2679\", \"80\u00d728=2240\"),\n    @(\"56\u00d765=3640\", \"54\u00d756=3024\"),\n    @(\"93\u00d767=6231\", \"82\u00d771=5822\"),\n    @(\"72\u00d748=3456\", \"54\u00d719=1026\"),\n    @(\"68\u00d758=3944\", \"55\u00d715=825\"),\n    @(\"14\u00d797=1358\", \"55\u00d760=3300\"),\n    @(\"17\u00d720=340\", \"50\u00d755=2750\"),\n    @(\"35\u00d737=1295\", \"68\u00d775=5100\"),\n    @(\"69\u00d733=2277\", \"57\u00d781=4617\"),\n    @(\"33\u00d777=2541\", \"66\u00d767=4422\"),\n    @(\"90\u00d799=8910\", \"16\u00d794=1504\"),\n    @(\"65\u00d768=4420\", \"83\u00d785=7055\"),\n    @(\"78\u00d762=4836\", \"51\u00d742=2142\"),\n    @(\"44\u00d786=3784\", \"40\u00d789=3560\"),\n    @(\"11\u00d724=264\", \"25\u00d721=525\"),\n    @(\"50\u00d730=1500\", \"33\u00d714=462\"),\n    @(\"60\u00d738=2280\", \"30\u00d796=2880\"),\n    @(\"68\u00d722=1496\", \"90\u00d7100=9000\"),\n    @(\"96\u00d733=3168\", \"71\u00d767=4757\"),\n    @(\"13\u00d729=377\", \"79\u00d745=3555\"),\n    @(\"22\u00d724=528\", \"28\u00d777=2156\"),\n    @(\"53\u00d737=1961\", \"19\u00d741=779\"),\n    @(\"38\u00d737=1406\", \"38\u00d755=2090\"),\n    @(\"95\u00d782=7790\", \"30\u00d763=1890\"),\n    @(\"76\u00d762=4712\", \"41\u00d7100=4100\"),\n    @(\"78\u00d713=1014\", \"59\u00d780=4720\"),\n    @(\"73\u00d789=6497\", \"71\u00d767=4757\"),\n    @(\"60\u00d793=5580\", \"30\u00d761=1830\"),\n    @(\"12\u00d777=924\", \"54\u00d793=5022\"),\n    @(\"37\u00d741=1517\", \"36\u00d789=3204\"),\n    @(\"88\u00d785=7480\", \"61\u00d799=6039\"),\n    @(\"28\u00d794=2632\", \"98\u00d713=1274\"),\n    @(\"64\u00d797=6208\", \"37\u00d738=1406\"),\n    @(\"72\u00d719=1368\", \"38\u00d780=3040\"),\n    @(\"41\u00d743=1763\", \"86\u00d720=1720\"),\n    @(\"50\u00d786=4300\", \"47\u00d795=4465\"),\n    @(\"50\u00d793=4650\", \"67\u00d790=6030\"),\n    @(\"19\u00d710=190\", \"82\u00d757=4674\"),\n    @(\"84\u00d791=7644\", \"36\u00d715=540\"),\n    @(\"89\u00d797=8633\", \"41\u00d777=3157\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n\"done\""}
